# Note: this runtime's Find.Execute (even when invoked on a sub-Range,
# e.g. Cell.Range.Find) searches/replaces across the whole document
# instead of being scoped to that Range. Several cells in this table
# share identical "before" text (e.g. "24+7=31" appears twice, mapped
# to two different results), so a scoped, position-based edit is
# required. Assigning directly to Range.Text is properly scoped to the
# specific paragraph/cell and preserves the existing run formatting
# (rFonts/sz), so that's used here instead of Find/Replace.

$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Paragraphs.Item(1).Range.Text = "2025-09-21 Sunday"

$t = $d.Tables.Item(1)

$t.Rows.Item(1).Cells.Item(1).Range.Text = "94-75=19"
$t.Rows.Item(1).Cells.Item(2).Range.Text = "32-6=26"
$t.Rows.Item(1).Cells.Item(3).Range.Text = "36+28=64"
$t.Rows.Item(1).Cells.Item(4).Range.Text = "73-44=29"
$t.Rows.Item(1).Cells.Item(5).Range.Text = "80-1=79"

$t.Rows.Item(2).Cells.Item(1).Range.Text = "34+58=92"
$t.Rows.Item(2).Cells.Item(2).Range.Text = "33-4=29"
$t.Rows.Item(2).Cells.Item(3).Range.Text = "38+39=77"
$t.Rows.Item(2).Cells.Item(4).Range.Text = "55+36=91"
$t.Rows.Item(2).Cells.Item(5).Range.Text = "30-21=9"

$t.Rows.Item(3).Cells.Item(1).Range.Text = "20-14=6"
$t.Rows.Item(3).Cells.Item(2).Range.Text = "27+15=42"
$t.Rows.Item(3).Cells.Item(3).Range.Text = "53-35=18"
$t.Rows.Item(3).Cells.Item(4).Range.Text = "31-14=17"
$t.Rows.Item(3).Cells.Item(5).Range.Text = "65-59=6"

$t.Rows.Item(4).Cells.Item(1).Range.Text = "77-38=39"
$t.Rows.Item(4).Cells.Item(2).Range.Text = "91-78=13"
$t.Rows.Item(4).Cells.Item(3).Range.Text = "37+45=82"
$t.Rows.Item(4).Cells.Item(4).Range.Text = "12-8=4"
$t.Rows.Item(4).Cells.Item(5).Range.Text = "17+45=62"

$t.Rows.Item(5).Cells.Item(1).Range.Text = "56+29=85"
$t.Rows.Item(5).Cells.Item(2).Range.Text = "5+28=33"
$t.Rows.Item(5).Cells.Item(3).Range.Text = "19+24=43"
$t.Rows.Item(5).Cells.Item(4).Range.Text = "9+87=96"
$t.Rows.Item(5).Cells.Item(5).Range.Text = "67-58=9"

$t.Rows.Item(6).Cells.Item(1).Range.Text = "14+67=81"
$t.Rows.Item(6).Cells.Item(2).Range.Text = "45+38=83"
$t.Rows.Item(6).Cells.Item(3).Range.Text = "96-69=27"
$t.Rows.Item(6).Cells.Item(4).Range.Text = "15+57=72"
$t.Rows.Item(6).Cells.Item(5).Range.Text = "7+68=75"

$t.Rows.Item(7).Cells.Item(1).Range.Text = "52+39=91"
$t.Rows.Item(7).Cells.Item(2).Range.Text = "33-28=5"
$t.Rows.Item(7).Cells.Item(3).Range.Text = "80-24=56"
$t.Rows.Item(7).Cells.Item(4).Range.Text = "18+43=61"
$t.Rows.Item(7).Cells.Item(5).Range.Text = "67+25=92"

$t.Rows.Item(8).Cells.Item(1).Range.Text = "6+65=71"
$t.Rows.Item(8).Cells.Item(2).Range.Text = "81-36=45"
$t.Rows.Item(8).Cells.Item(3).Range.Text = "24+69=93"
$t.Rows.Item(8).Cells.Item(4).Range.Text = "8+83=91"
$t.Rows.Item(8).Cells.Item(5).Range.Text = "71-26=45"

$t.Rows.Item(9).Cells.Item(1).Range.Text = "71-63=8"
$t.Rows.Item(9).Cells.Item(2).Range.Text = "57+4=61"
$t.Rows.Item(9).Cells.Item(3).Range.Text = "55+8=63"
$t.Rows.Item(9).Cells.Item(4).Range.Text = "71-36=35"
$t.Rows.Item(9).Cells.Item(5).Range.Text = "88+9=97"

$t.Rows.Item(10).Cells.Item(1).Range.Text = "82-13=69"
$t.Rows.Item(10).Cells.Item(2).Range.Text = "16-9=7"
$t.Rows.Item(10).Cells.Item(3).Range.Text = "82-67=15"
$t.Rows.Item(10).Cells.Item(4).Range.Text = "79+18=97"
$t.Rows.Item(10).Cells.Item(5).Range.Text = "82-28=54"

$t.Rows.Item(11).Cells.Item(1).Range.Text = "53-4=49"
$t.Rows.Item(11).Cells.Item(2).Range.Text = "22+9=31"
$t.Rows.Item(11).Cells.Item(3).Range.Text = "72-27=45"
$t.Rows.Item(11).Cells.Item(4).Range.Text = "74+18=92"
$t.Rows.Item(11).Cells.Item(5).Range.Text = "61-26=35"

$t.Rows.Item(12).Cells.Item(1).Range.Text = "59+34=93"
$t.Rows.Item(12).Cells.Item(2).Range.Text = "43-28=15"
$t.Rows.Item(12).Cells.Item(3).Range.Text = "43-18=25"
$t.Rows.Item(12).Cells.Item(4).Range.Text = "56+18=74"
$t.Rows.Item(12).Cells.Item(5).Range.Text = "30-8=22"

$t.Rows.Item(13).Cells.Item(1).Range.Text = "11-4=7"
$t.Rows.Item(13).Cells.Item(2).Range.Text = "95-78=17"
$t.Rows.Item(13).Cells.Item(3).Range.Text = "82-66=16"
$t.Rows.Item(13).Cells.Item(4).Range.Text = "27+56=83"
$t.Rows.Item(13).Cells.Item(5).Range.Text = "43-28=15"

$t.Rows.Item(14).Cells.Item(1).Range.Text = "33-18=15"
$t.Rows.Item(14).Cells.Item(2).Range.Text = "14+8=22"
$t.Rows.Item(14).Cells.Item(3).Range.Text = "59+8=67"
$t.Rows.Item(14).Cells.Item(4).Range.Text = "41-26=15"
$t.Rows.Item(14).Cells.Item(5).Range.Text = "66+6=72"

$t.Rows.Item(15).Cells.Item(1).Range.Text = "39+57=96"
$t.Rows.Item(15).Cells.Item(2).Range.Text = "95-86=9"
$t.Rows.Item(15).Cells.Item(3).Range.Text = "70-17=53"
$t.Rows.Item(15).Cells.Item(4).Range.Text = "95-79=16"
$t.Rows.Item(15).Cells.Item(5).Range.Text = "83-29=54"

$t.Rows.Item(16).Cells.Item(1).Range.Text = "92-89=3"
$t.Rows.Item(16).Cells.Item(2).Range.Text = "9+88=97"
$t.Rows.Item(16).Cells.Item(3).Range.Text = "82-74=8"
$t.Rows.Item(16).Cells.Item(4).Range.Text = "66-49=17"
$t.Rows.Item(16).Cells.Item(5).Range.Text = "19+76=95"

$t.Rows.Item(17).Cells.Item(1).Range.Text = "90-57=33"
$t.Rows.Item(17).Cells.Item(2).Range.Text = "39+44=83"
$t.Rows.Item(17).Cells.Item(3).Range.Text = "19+69=88"
$t.Rows.Item(17).Cells.Item(4).Range.Text = "27+8=35"
$t.Rows.Item(17).Cells.Item(5).Range.Text = "7+35=42"

$t.Rows.Item(18).Cells.Item(1).Range.Text = "90-59=31"
$t.Rows.Item(18).Cells.Item(2).Range.Text = "18+14=32"
$t.Rows.Item(18).Cells.Item(3).Range.Text = "91-83=8"
$t.Rows.Item(18).Cells.Item(4).Range.Text = "46+16=62"
$t.Rows.Item(18).Cells.Item(5).Range.Text = "81-68=13"

$t.Rows.Item(19).Cells.Item(1).Range.Text = "15+67=82"
$t.Rows.Item(19).Cells.Item(2).Range.Text = "62-24=38"
$t.Rows.Item(19).Cells.Item(3).Range.Text = "7+37=44"
$t.Rows.Item(19).Cells.Item(4).Range.Text = "49+45=94"
$t.Rows.Item(19).Cells.Item(5).Range.Text = "60-12=48"

$t.Rows.Item(20).Cells.Item(1).Range.Text = "84-36=48"
$t.Rows.Item(20).Cells.Item(2).Range.Text = "13+38=51"
$t.Rows.Item(20).Cells.Item(3).Range.Text = "35+39=74"
$t.Rows.Item(20).Cells.Item(4).Range.Text = "46+27=73"
$t.Rows.Item(20).Cells.Item(5).Range.Text = "17+18=35"

Write-Host "All replacements complete"